# Insert a new weekly record at the top of the data block (row 109),
# pushing all existing records (old rows 109-231) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("109:109").Insert()

$ws.Range("A109").Value = 5
$ws.Range("B109").Value = 'Macroferia Regional de Talca'
$ws.Range("C109").Value = 'Maule'
$ws.Range("D109").Value = 44994
$ws.Range("E109").Value = 7
$ws.Range("F109").Value = 100112031
$ws.Range("G109").Value = 'Poroto verde'
$ws.Range("H109").Value = 'Sin especificar'
$ws.Range("I109").Value = 'Primera'
$ws.Range("J109").Value = 200
$ws.Range("K109").Value = 25000
$ws.Range("L109").Value = 25000
$ws.Range("M109").Value = 25000
$ws.Range("N109").Value = '$/saco 25 kilos'
$ws.Range("O109").Value = 'Región del Maule'
$ws.Range("P109").Value = 1000
$ws.Range("Q109").Value = 25
$ws.Range("R109").Value = 'Hortaliza'
